$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02453504007225149
$ws.Range("D2").Value = 0.02595913159737506
$ws.Range("E2").Value = 0.1083956125688772
$ws.Range("F2").Value = 0.4011173264326615
$ws.Range("G2").Value = 0.2494694137112177
$ws.Range("H2").Value = 0.4299996397747279
$ws.Range("I2").Value = 0.294091890236551
$ws.Range("K2").Value = 1.187768794381611
$ws.Range("M2").Value = 0.3730441353413028
$ws.Range("O2").Value = 1.276956331357368

$ws.Range("C3").Value = 0.02138906953582875
$ws.Range("D3").Value = 0.02290482479835276
$ws.Range("E3").Value = 0.1051323248934608
$ws.Range("F3").Value = 0.402784158121726
$ws.Range("G3").Value = 0.2521603659494431
$ws.Range("H3").Value = 0.4354858015413043
$ws.Range("I3").Value = 0.2960398840774836
$ws.Range("K3").Value = 1.036841494584053
$ws.Range("M3").Value = 0.3303842439551445
$ws.Range("O3").Value = 1.293880245128435

$ws.Range("C4").Value = 0.0194513134007579
$ws.Range("D4").Value = 0.02102014890999726
$ws.Range("E4").Value = 0.1032786869389888
$ws.Range("F4").Value = 0.4042073267939728
$ws.Range("G4").Value = 0.2541391123286942
$ws.Range("H4").Value = 0.4391444145723113
$ws.Range("I4").Value = 0.2975333082031355
$ws.Range("K4").Value = 0.9437758093563389
$ws.Range("M4").Value = 0.3041870526631314
$ws.Range("O4").Value = 1.305562003549809

$ws.Range("C5").Value = 0.0186601686363872
$ws.Range("D5").Value = 0.02024984000026819
$ws.Range("E5").Value = 0.1025606987450729
$ws.Range("F5").Value = 0.4048875532774616
$ws.Range("G5").Value = 0.2550272263412197
$ws.Range("H5").Value = 0.4407081995181557
$ws.Range("I5").Value = 0.2982164714525162
$ws.Range("K5").Value = 0.9057540600553011
$ws.Range("M5").Value = 0.293510857562616
$ws.Range("O5").Value = 1.310646015819373

$ws.Range("C6").Value = 0.01852871077925045
$ws.Range("D6").Value = 0.02012179422303007
$ws.Range("E6").Value = 0.1024437252344832
$ws.Range("F6").Value = 0.4050065535628775
$ws.Range("G6").Value = 0.2551796250344438
$ws.Range("H6").Value = 0.4409722638257847
$ws.Range("I6").Value = 0.2983344091048608
$ws.Range("K6").Value = 0.8994348177621987
$ws.Range("M6").Value = 0.29173805576945
$ws.Range("O6").Value = 1.311509730945886

$ws.Range("C7").Value = 0.01944064973363169
$ws.Range("D7").Value = 0.0210097694457545
$ws.Range("E7").Value = 0.1032688529852201
$ws.Range("F7").Value = 0.4042160948826279
$ws.Range("G7").Value = 0.2541507591587973
$ws.Range("H7").Value = 0.4391652093725753
$ws.Range("I7").Value = 0.297542219859448
$ws.Range("K7").Value = 0.943263422689455
$ws.Range("M7").Value = 0.3040430718453493
$ws.Range("O7").Value = 1.305629259244284

$ws.Range("C8").Value = 0.02345160506169464
$ws.Range("D8").Value = 0.02490797291308411
$ws.Range("E8").Value = 0.1072390839358945
$ws.Range("F8").Value = 0.4016089410628041
$ws.Range("G8").Value = 0.250329286590059
$ws.Range("H8").Value = 0.4318310162744794
$ws.Range("I8").Value = 0.2947017295349674
$ws.Range("K8").Value = 1.135812805101523
$ws.Range("M8").Value = 0.3583359925045428
$ws.Range("O8").Value = 1.282523343268792

$ws.Range("C9").Value = 0.03126693897377208
$ws.Range("D9").Value = 0.03247633163179842
$ws.Range("E9").Value = 0.1162311213030449
$ws.Range("F9").Value = 0.3996793260392124
$ws.Range("G9").Value = 0.2454413810785638
$ws.Range("H9").Value = 0.4197537116403467
$ws.Range("I9").Value = 0.2914996847947684
$ws.Range("K9").Value = 1.510158282307771
$ws.Range("M9").Value = 0.4647642679878174
$ws.Range("O9").Value = 1.24749143031687

$ws.Range("C10").Value = 0.03697670434927147
$ws.Range("D10").Value = 0.03798822899842946
$ws.Range("E10").Value = 0.123595882254655
$ws.Range("F10").Value = 0.4002184817530861
$ws.Range("G10").Value = 0.2434604149440531
$ws.Range("H10").Value = 0.4122902402783666
$ws.Range("I10").Value = 0.2906034238119872
$ws.Range("K10").Value = 1.783100990444893
$ws.Range("M10").Value = 0.542928115475334
$ws.Range("O10").Value = 1.228075076353932

$ws.Range("C11").Value = 0.03956693825264779
$ws.Range("D11").Value = 0.04048474807760272
$ws.Range("E11").Value = 0.1271159174210226
$ws.Range("F11").Value = 0.4008922072418315
$ws.Range("G11").Value = 0.2429134581349217
$ws.Range("H11").Value = 0.4092019933884927
$ws.Range("I11").Value = 0.2905146210257001
$ws.Range("K11").Value = 1.906794490720188
$ws.Range("M11").Value = 0.5784801611955004
$ws.Range("O11").Value = 1.220626616020596

$ws.Range("C12").Value = 0.04054672404106441
$ws.Range("D12").Value = 0.04142850367013295
$ws.Range("E12").Value = 0.1284736603205303
$ws.Range("F12").Value = 0.4012092107631773
$ws.Range("G12").Value = 0.2427576386998922
$ws.Range("H12").Value = 0.4080767735279451
$ws.Range("I12").Value = 0.2905270625453795
$ws.Range("K12").Value = 1.953564125171624
$ws.Range("M12").Value = 0.5919418548917008
$ws.Range("O12").Value = 1.218006072339676

$ws.Range("C13").Value = 0.04033575843199344
$ws.Range("D13").Value = 0.04122532208874929
$ws.Range("E13").Value = 0.1281801380232466
$ws.Range("F13").Value = 0.4011381821154174
$ws.Range("G13").Value = 0.2427889097217957
$ws.Range("H13").Value = 0.4083171409784399
$ws.Range("I13").Value = 0.2905223307708127
$ws.Range("K13").Value = 1.943494624843765
$ws.Range("M13").Value = 0.5890426927143722
$ws.Range("O13").Value = 1.218561541373575

$ws.Range("C14").Value = 0.03964756773603995
$ws.Range("D14").Value = 0.04056242426105428
$ws.Range("E14").Value = 0.127227120961507
$ws.Range("F14").Value = 0.4009170456627089
$ws.Range("G14").Value = 0.2428996085197426
$ws.Range("H14").Value = 0.4091085338399481
$ws.Range("I14").Value = 0.2905147201030545
$ws.Range("K14").Value = 1.91064368428556
$ws.Range("M14").Value = 0.5795876868398153
$ws.Range("O14").Value = 1.220407007835377

$ws.Range("C15").Value = 0.03922588902031521
$ws.Range("D15").Value = 0.04015616722911375
$ws.Range("E15").Value = 0.1266466088080449
$ws.Range("F15").Value = 0.4007896594413083
$ws.Range("G15").Value = 0.2429741065654412
$ws.Range("H15").Value = 0.4095990480295981
$ws.Range("I15").Value = 0.2905160641316229
$ws.Range("K15").Value = 1.890512279816164
$ws.Range("M15").Value = 0.5737960694848141
$ws.Range("O15").Value = 1.221563487971352

$ws.Range("C16").Value = 0.03680727806273865
$ws.Range("D16").Value = 0.03782485182348694
$ws.Range("E16").Value = 0.1233692874331069
$ws.Range("F16").Value = 0.4001830967397595
$ws.Range("G16").Value = 0.2435033187066011
$ws.Range("H16").Value = 0.4124982488623132
$ws.Range("I16").Value = 0.2906156625909517
$ws.Range("K16").Value = 1.775007633464497
$ws.Range("M16").Value = 0.540604572235452
$ws.Range("O16").Value = 1.228589787822301

$ws.Range("C17").Value = 0.03532166803694281
$ws.Range("D17").Value = 0.03639183912289923
$ws.Range("E17").Value = 0.1214025015255373
$ws.Range("F17").Value = 0.3999209193767399
$ws.Range("G17").Value = 0.2439189538100877
$ws.Range("H17").Value = 0.4143554943009065
$ws.Range("I17").Value = 0.2907585938196036
$ws.Range("K17").Value = 1.70402695281831
$ws.Range("M17").Value = 0.5202411021548272
$ws.Range("O17").Value = 1.233255411479433

$ws.Range("C18").Value = 0.03446651276718171
$ws.Range("D18").Value = 0.03556658841002047
$ws.Range("E18").Value = 0.120287215606993
$ws.Range("F18").Value = 0.3998104348895737
$ws.Range("G18").Value = 0.2441913327083256
$ws.Range("H18").Value = 0.4154526160887784
$ws.Range("I18").Value = 0.2908708090126808
$ws.Range("K18").Value = 1.663156721650296
$ws.Range("M18").Value = 0.5085281160198463
$ws.Range("O18").Value = 1.236069177464572

$ws.Range("C19").Value = 0.03417685802006076
$ws.Range("D19").Value = 0.03528699936131829
$ws.Range("E19").Value = 0.1199123283874712
$ws.Range("F19").Value = 0.3997799423903956
$ws.Range("G19").Value = 0.2442892664163807
$ws.Range("H19").Value = 0.4158290404558045
$ws.Range("I19").Value = 0.2909139500938771
$ws.Range("K19").Value = 1.649311291672973
$ws.Range("M19").Value = 0.5045622328857036
$ws.Range("O19").Value = 1.237044205914941

$ws.Range("C20").Value = 0.03547988363133925
$ws.Range("D20").Value = 0.03654449173413354
$ws.Range("E20").Value = 0.1216102150457701
$ws.Range("F20").Value = 0.3999446544250063
$ws.Range("G20").Value = 0.2438712579265996
$ws.Range("H20").Value = 0.4141547973319604
$ws.Range("I20").Value = 0.2907402716203933
$ws.Range("K20").Value = 1.71158754452432
$ws.Range("M20").Value = 0.5224088782311753
$ws.Range("O20").Value = 1.232745263133125

$ws.Range("C21").Value = 0.03984973571370176
$ws.Range("D21").Value = 0.04075717790924216
$ws.Range("E21").Value = 0.1275063696526288
$ws.Range("F21").Value = 0.4009803173302871
$ws.Range("G21").Value = 0.2428656983724551
$ws.Range("H21").Value = 0.4088748813526024
$ws.Range("I21").Value = 0.2905157035681611
$ws.Range("K21").Value = 1.920294735667994
$ws.Range("M21").Value = 0.5823648838380109
$ws.Range("O21").Value = 1.219859513077367

$ws.Range("C22").Value = 0.04269937145713243
$ws.Range("D22").Value = 0.04350093056825699
$ws.Range("E22").Value = 0.1315044674942385
$ws.Range("F22").Value = 0.4020179788831513
$ws.Range("G22").Value = 0.2425076994561692
$ws.Range("H22").Value = 0.4056820149454268
$ws.Range("I22").Value = 0.290637558384212
$ws.Range("K22").Value = 2.056285226371017
$ws.Range("M22").Value = 0.6215430545728822
$ws.Range("O22").Value = 1.212604293981812

$ws.Range("C23").Value = 0.04117906120801251
$ws.Range("D23").Value = 0.04203742500693863
$ws.Range("E23").Value = 0.1293572525543496
$ws.Range("F23").Value = 0.4014310591848727
$ws.Range("G23").Value = 0.2426712728896234
$ws.Range("H23").Value = 0.4073624813374863
$ws.Range("I23").Value = 0.290547874108519
$ws.Range("K23").Value = 1.98374316379784
$ws.Range("M23").Value = 0.6006336454195633
$ws.Range("O23").Value = 1.216369493083036

$ws.Range("C24").Value = 0.03540835765370787
$ws.Range("D24").Value = 0.03647548182216553
$ws.Range("E24").Value = 0.1215162596814281
$ws.Range("F24").Value = 0.3999337984703573
$ws.Range("G24").Value = 0.2438927171477445
$ws.Range("H24").Value = 0.414245440962155
$ws.Range("I24").Value = 0.2907484615233678
$ws.Range("K24").Value = 1.708169595454592
$ws.Range("M24").Value = 0.5214288445732365
$ws.Range("O24").Value = 1.232975491856493

$ws.Range("C25").Value = 0.02915821729381207
$ws.Range("D25").Value = 0.03043726019056692
$ws.Range("E25").Value = 0.1136670203445149
$ws.Range("F25").Value = 0.399858872300122
$ws.Range("G25").Value = 0.2464824704359856
$ws.Range("H25").Value = 0.4227737053531726
$ws.Range("I25").Value = 0.2921110887293743
$ws.Range("K25").Value = 1.409246793406453
$ws.Range("M25").Value = 0.4359773767467487
$ws.Range("O25").Value = 1.255862337557687
